$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New line entries inserted in the string table shift the "extr" labels down
# by two positions for rows 8-15. Row 8/9 become the new line7/line8 entries,
# and rows 10-15 take on the name that used to belong two rows above them.

# Row 8 (was extr1, now line7): B,C,D,E updated
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# Row 9 (was extr2, now line8): B,C updated
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16

# Row 10 (was extr3, now extr1): B,C,D,E updated
$ws.Cells.Item(10, 2).Value = "extr1"
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

# Row 11 (was extr4, now extr2): B,C,D,E updated
$ws.Cells.Item(11, 2).Value = "extr2"
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true

# Row 12 (was extr5, now extr3): B,C updated
$ws.Cells.Item(12, 2).Value = "extr3"
$ws.Cells.Item(12, 3).Value = 10

# Row 13 (was extr6, now extr4): B,D updated
$ws.Cells.Item(13, 2).Value = "extr4"
$ws.Cells.Item(13, 4).Value = 8

# Row 14 (was extr7, now extr5): B,C,D updated
$ws.Cells.Item(14, 2).Value = "extr5"
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11

# Row 15 (was extr8, now extr6): B,C,D,E updated
$ws.Cells.Item(15, 2).Value = "extr6"
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $true

# New row 16: A=14, B="extr7", C=5, D=7, E=false
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "extr7"
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $false

# New row 17: A=15, B="extr8", C=8, D=5, E=false
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "extr8"
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $false

# Copy formatting (bold, centered, bordered) from the existing column-A data
# style onto the two new rows, matching the rest of column A (A2:A15).
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16:A17").PasteSpecial(-4122) | Out-Null
